{"js": "// Add a \"Sep 11, 2024\" paragraph right after the \"Movies: Setup\" paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Movies: Setup\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph with text \"Movies: Setup\"');\n}\n\ntarget.insertParagraph(\"Sep 11, 2024\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Add a \"Sep 11, 2024\" paragraph right after the \"Movies: Setup\" paragraph.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`a\", \"`n\") -eq \"Movies: Setup\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    $r = $d.Content\n    $r.Find.Execute(\"Movies: Setup\") | Out-Null\n    $target = $r.Paragraphs(1)\n}\n\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.Text = \"Sep 11, 2024\"\n"}
